$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Extend the centered-alignment formatting (style used by A5:J22) down
#        into the newly used rows, matching exactly the ranges that receive
#        cells in the target file (A23:E23, A24:F24, A25:E28). Copy/PasteSpecial
#        (formats only) reuses the existing style index instead of minting a
#        new, unused cellXf entry. ---
$ws.Range("A21:E21").Copy()
$ws.Range("A23:E23").PasteSpecial(-4122)
$ws.Range("A25:E25").PasteSpecial(-4122)
$ws.Range("A26:E26").PasteSpecial(-4122)
$ws.Range("A27:E27").PasteSpecial(-4122)
$ws.Range("A28:E28").PasteSpecial(-4122)

$ws.Range("A21:F21").Copy()
$ws.Range("A24:F24").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- 2. Fill in the new component rows. Values are entered in the same
#        order the original author typed them (this controls the order new
#        shared strings are appended in, matching the source diff). ---

# Row 21: Potenciometro 5K
$ws.Range("B21").Value = 1
$ws.Range("E21").Value = "5KΩ"
$ws.Range("D21").Value = "Potenciômetro"
$ws.Range("C21").Value = "PTH"

# Row 22: LDR 5mm
$ws.Range("B22").Value = 1
$ws.Range("C22").Value = "PTH 5mm"
$ws.Range("D22").Value = "LDR 5mm"
$ws.Range("E22").Value = "10KΩ"

# Row 23: NRF24l01
$ws.Range("B23").Value = 1
$ws.Range("C23").Value = "PTH"
$ws.Range("D23").Value = "NRF24l01"
$ws.Range("E23").Value = "-"

# Row 24: Bluetooth
$ws.Range("B24").Value = 1
$ws.Range("C24").Value = "PTH"
$ws.Range("D24").Value = "Bluetooth"
$ws.Range("E24").Value = "-"

# "OU" (OR) annotations - written last, matching original edit order
$ws.Range("F21").Value = "OU"
$ws.Range("F24").Value = "OU"

# --- 3. Restore the on-screen selection to match the saved view state. ---
$ws.Range("G20").Select()

$wb.Save()
